# Harmonize Excel configuration file:
# - Append trailing slashes to folder path values
# - Rename several "Property" keys to their harmonized names
# - Move the active selection to C14

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2: modelFolder -> add trailing slash to folder value
$ws.Range("B2").Value = "Models/Simulations/"

# Row 3: paramsFolder -> add trailing slash to folder value
$ws.Range("B3").Value = "Parameters/"

# Row 4: paramsFile -> modelParamsFile
$ws.Range("A4").Value = "modelParamsFile"

# Row 6: populationParamsFile -> populationsFile, PopulationParameters.xlsx -> Populations.xlsx
$ws.Range("A6").Value = "populationsFile"
$ws.Range("B6").Value = "Populations.xlsx"

# Row 7: scenarioDefinitionFile -> scenariosFile
$ws.Range("A7").Value = "scenariosFile"

# Row 8: scenarioApplicationsFile -> applicationsFile, ApplicationParameters.xlsx -> Applications.xlsx
$ws.Range("A8").Value = "applicationsFile"
$ws.Range("B8").Value = "Applications.xlsx"

# Row 10: dataFolder -> add trailing slash to folder value
$ws.Range("B10").Value = "Data/"

# Row 14: outputFolder -> add trailing slash to folder value
$ws.Range("B14").Value = "Results/"

# Update the selected/active cell to C14
$ws.Range("C14").Select()
